# removing thread sleeps from profile testcases
# Add a new F-suite test case (TestCase_F19 / OPQA-1012) to the "Test Cases" sheet
# and flip the Results of row 2 (TestCase_F1) from SKIP to PASS.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")  # "Test Cases" is the active/tab-selected sheet

# New row 20: TestCase_F19
# (values set in this order so new shared strings are interned in the
#  same sequence as the target workbook: description, jira id, tcid, PASS)
$ws.Range("C20").Value = "Verify that follower of the article is able to start conversation from home page when some one commented on the article he is following."
$ws.Range("B20").Value = "OPQA-1012"
$ws.Range("A20").Value = "TestCase_F19"
$ws.Range("D20").Value = "Y"
$ws.Range("E20").Value = "PASS"

# Row 2 (TestCase_F1) Results: SKIP -> PASS (reuses the shared "PASS" string)
$ws.Range("E2").Value = "PASS"

# Copy formatting from existing rows so the new row matches the sheet's
# established look: A/D/E (TCID / Runmode / Results) take the shaded
# "s=6" style used throughout the column, while B/C (Jira id / Description)
# take the plain bordered "s=2" style used by the Description column.
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E19").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null

# Update selection to match target workbook state (the sheet also scrolls
# so column B is the leftmost visible column, but horizontal scroll
# position isn't part of the exposed COM object model here).
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("D13").Select() | Out-Null

$wb.Save()
